# Update cryptos list data (prices and volume% changes) per daily refresh
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "29.263.36"
$ws.Range("E2").Value = "  +0.50%  "
$ws.Range("D3").Value = "1.891.12"
$ws.Range("E3").Value = "  -0.42%  "
$ws.Range("D4").Value = "'1.000"
$ws.Range("E4").Value = "  -0.04%  "
$ws.Range("D5").Value = "322.73"
$ws.Range("E5").Value = "  -2.95%  "
$ws.Range("D6").Value = "0.9996"
$ws.Range("E6").Value = "  -0.10%  "
$ws.Range("D7").Value = "0.4712"
$ws.Range("E7").Value = "  +2.38%  "
$ws.Range("D8").Value = "0.4037"
$ws.Range("E8").Value = "  -1.92%  "
$ws.Range("D9").Value = "47.37"
$ws.Range("E9").Value = "  -0.97%  "
$ws.Range("D10").Value = "0.08015"
$ws.Range("E10").Value = "  +0.04%  "
$ws.Range("D11").Value = "0.9946"
$ws.Range("E11").Value = "  -1.55%  "
$ws.Range("D12").Value = "22.83"
$ws.Range("E12").Value = "  +3.29%  "
$ws.Range("D13").Value = "1.884.35"
$ws.Range("E13").Value = "  -0.64%  "
$ws.Range("D14").Value = "'5.930"
$ws.Range("E14").Value = "  -0.04%  "
$ws.Range("D15").Value = "7.039"
$ws.Range("E15").Value = "  -0.96%  "
$ws.Range("D16").Value = "89.39"
$ws.Range("E16").Value = "  +0.31%  "
$ws.Range("E17").Value = "  -0.02%  "
$ws.Range("D18").Value = "0.06628"
$ws.Range("E18").Value = "  +0.94%  "
$ws.Range("D19").Value = "'0.00001023"
$ws.Range("E19").Value = "  -0.64%  "
$ws.Range("D20").Value = "17.45"
$ws.Range("E20").Value = "  -0.93%  "
$ws.Range("E21").Value = "  +0.06%  "
$ws.Range("D22").Value = "29.261.56"
$ws.Range("E22").Value = "  +0.55%  "
$ws.Range("D23").Value = "5.492"
$ws.Range("E23").Value = "  -0.07%  "
$ws.Range("D24").Value = "11.68"
$ws.Range("E24").Value = "  +2.50%  "
$ws.Range("D25").Value = "2.174"
$ws.Range("E25").Value = "  -0.84%  "
$ws.Range("D26").Value = "2.139.48"
$ws.Range("E26").Value = "  +0.89%  "
$ws.Range("D27").Value = "155.12"
$ws.Range("E27").Value = "  -1.20%  "
$ws.Range("D28").Value = "19.64"
$ws.Range("E28").Value = "  -0.53%  "
$ws.Range("D29").Value = "5.994"
$ws.Range("E29").Value = "  +6.81%  "
$ws.Range("D30").Value = "2.085"
$ws.Range("E30").Value = "  -1.87%  "
$ws.Range("D31").Value = "117.02"
$ws.Range("E31").Value = "  -0.09%  "
$ws.Range("E32").Value = "  -2.12%  "
$ws.Range("D33").Value = "0.09414"
$ws.Range("E33").Value = "  +0.18%  "
$ws.Range("D34").Value = "3.537"
$ws.Range("E34").Value = "  -0.31%  "
$ws.Range("D35").Value = "1.381"
$ws.Range("E35").Value = "  -2.62%  "
$ws.Range("D36").Value = "5.353"
$ws.Range("E36").Value = "  +0.01%  "
$ws.Range("B37").Value = "Hedera"
$ws.Range("C37").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D37").Value = "0.06047"
$ws.Range("E37").Value = "  -0.70%  "
$ws.Range("B38").Value = "VeChain"
$ws.Range("C38").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D38").Value = "0.02241"
$ws.Range("E38").Value = "  -0.02%  "
$ws.Range("D39").Value = "1.168"
$ws.Range("E39").Value = "  -0.77%  "
$ws.Range("D40").Value = "8.006"
$ws.Range("E40").Value = "  -4.95%  "
$ws.Range("D41").Value = "0.5817"
$ws.Range("E41").Value = "  -0.41%  "
$ws.Range("D42").Value = "0.1831"
$ws.Range("E42").Value = "  +0.26%  "
$ws.Range("D43").Value = "10.05"
$ws.Range("E43").Value = "  -0.65%  "
$ws.Range("B44").Value = "WEMIXToken"
$ws.Range("C44").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D44").Value = "1.272"
$ws.Range("E44").Value = "  +1.41%  "
$ws.Range("B45").Value = "RenderToken"
$ws.Range("C45").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D45").Value = "2.383"
$ws.Range("E45").Value = "  +1.00%  "
$ws.Range("D46").Value = "'0.07710"
$ws.Range("E46").Value = "  +2.75%  "
$ws.Range("D47").Value = "12.15"
$ws.Range("E47").Value = "  +0.56%  "
$ws.Range("D48").Value = "0.5475"
$ws.Range("E48").Value = "  -1.25%  "
$ws.Range("D49").Value = "1.904"
$ws.Range("E49").Value = "  -1.00%  "
$ws.Range("D50").Value = "113.36"
$ws.Range("E50").Value = "  +0.85%  "
$ws.Range("D51").Value = "0.2956"
$ws.Range("E51").Value = "  +1.12%  "
